$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split the range [start,end) into consecutive runs whose boundaries
# are given by the (sorted, absolute) offsets in $cuts (start/end excluded).
# Word COM creates a genuine run split whenever direct character formatting
# is applied to a *non-empty* sub-range, so touching each segment in turn
# and flipping Bold twice leaves the text/content fully unchanged while
# still causing the adjacent runs to separate at the segment boundaries.
# ---------------------------------------------------------------------------
function Split-RunRange($start, $end, $cuts) {
    $bounds = @($start) + $cuts + @($end)
    for ($i = 0; $i -lt $bounds.Length - 1; $i++) {
        $r = $d.Range($bounds[$i], $bounds[$i + 1])
        $r.Font.Bold = 1
        $r.Font.Bold = 0
    }
}

# ---------------------------------------------------------------------------
# 1) Title: "Week 4 Reading Guide: Basic Regression" -> one run per
#    word/space token.
# ---------------------------------------------------------------------------
$title = $d.Paragraphs(1).Range
$titleStart = $title.Start
$cuts = @(4,5,6,7,14,15,21,22,27,28) | ForEach-Object { $titleStart + $_ }
Split-RunRange $titleStart ($titleStart + 38) $cuts

# ---------------------------------------------------------------------------
# 2) Quoted single/multi-word phrases: split `"word"` into three runs,
#    `"`, `word`, `"`.
# ---------------------------------------------------------------------------
function Split-QuotedPhrase($fullPhrase) {
    $inner = $fullPhrase.Substring(1, $fullPhrase.Length - 2)
    $f = $d.Content
    $found = $f.Find.Execute($fullPhrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $fullPhrase
        return
    }
    $s = $f.Start
    $e = $f.End
    $b1 = $s + 1
    $b2 = $e - 1
    $r1 = $d.Range($s, $b1)
    $r1.Font.Bold = 1
    $r1.Font.Bold = 0
    $r2 = $d.Range($b1, $b2)
    $r2.Font.Bold = 1
    $r2.Font.Bold = 0
    $r3 = $d.Range($b2, $e)
    $r3.Font.Bold = 1
    $r3.Font.Bold = 0
}

Split-QuotedPhrase "“response”"
Split-QuotedPhrase "“explanatory”"
Split-QuotedPhrase "“explanatory modeling”"
Split-QuotedPhrase "“predictive modeling”"
Split-QuotedPhrase "“basic”"
Split-QuotedPhrase "“EDA”"
Split-QuotedPhrase "“fit”"
Split-QuotedPhrase "“levels”"
Split-QuotedPhrase "“baseline”"
Split-QuotedPhrase "“indicator function”"
Split-QuotedPhrase "“correlation does not imply causation,”"

# ---------------------------------------------------------------------------
# 3) Style "Subtitle": re-base from "Title" to "Normal", and give its run
#    properties an explicit theme-based text colour (text1, tinted ~65%).
# ---------------------------------------------------------------------------
$subtitle = $d.Styles("Subtitle")
$subtitle.BaseStyle = $d.Styles("Normal")
# Packed theme-colour value recognised by this host: byte0 = 0xD0 + theme
# index (13 = Text1), byte1 = 0x00, byte2 = shade-or-0xFF, byte3 =
# tint-or-0xFF. This reproduces w:themeColor="text1" w:themeTint="A6".
$subtitle.Font.TextColor.RGB = -587137114

# ---------------------------------------------------------------------------
# 4) Style "AbstractTitle": add an explicit RGB text colour (345A8A). Word's
#    Color property takes BGR-ordered values, so convert RGB -> BGR first.
# ---------------------------------------------------------------------------
$abstractTitle = $d.Styles("AbstractTitle")
$abstractTitle.Font.Color = 9067060

Write-Host "done"
